$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on price cells whose new values would
# otherwise be auto-parsed as numbers by Excel, so they stay text
# (matching the original inline-string "Price" column formatting).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "52.195.62"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "2.893.64"
$ws.Range("E3").Value = "  +3.49%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "352.61"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "112.07"
$ws.Range("E6").Value = "  +2.75%  "
$ws.Range("D7").Value = "0.563"
$ws.Range("E7").Value = "  +1.77%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").Value = "40.11"
$ws.Range("E10").Value = "  +0.99%  "
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "0.0858"
$ws.Range("E12").Value = "  +2.80%  "
$ws.Range("D13").Value = "20.00"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "3.346.72"
$ws.Range("E15").Value = "  +3.62%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.908.07"
$ws.Range("E16").Value = "  +3.65%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "0.992"
$ws.Range("E17").Value = "  +6.15%  "
$ws.Range("D18").Value = "52.172.36"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "7.74"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "3.35"
$ws.Range("E20").Value = "  +6.58%  "
$ws.Range("D21").Value = "14.48"
$ws.Range("E21").Value = "  +8.17%  "
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("D23").Value = "71.04"
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("D24").Value = "270.50"
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("D26").Value = "26.47"
$ws.Range("E26").Value = "  +2.07%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "0.164"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").Value = "38.89"
$ws.Range("E29").Value = "  +4.08%  "
$ws.Range("D30").Value = "10.50"
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("E31").Value = "  +0.89%  "
$ws.Range("D32").Value = "6.48"
$ws.Range("E32").Value = "  +3.48%  "
$ws.Range("D33").Value = "6.14"
$ws.Range("E33").Value = "  +7.90%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "53.35"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.0950"
$ws.Range("E35").Value = "  +11.22%  "
$ws.Range("E36").Value = "  +3.43%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  +6.16%  "
$ws.Range("D39").Value = "18.66"
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("E41").Value = "  +5.81%  "
$ws.Range("E42").Value = "  +2.59%  "
$ws.Range("D43").Value = "22.78"
$ws.Range("E43").Value = "  +3.77%  "
$ws.Range("D44").Value = "121.80"
$ws.Range("E44").Value = "  +1.46%  "
$ws.Range("D45").Value = "2.20"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").Value = "3.59"
$ws.Range("E46").Value = "  +6.53%  "
$ws.Range("D47").Value = "2.207.04"
$ws.Range("E47").Value = "  +3.33%  "
$ws.Range("E48").Value = "  +6.23%  "
$ws.Range("D49").Value = "0.268"
$ws.Range("E49").Value = "  +18.79%  "
$ws.Range("D50").Value = "0.952"
$ws.Range("E50").Value = "  +4.74%  "
$ws.Range("D51").Value = "5.51"
$ws.Range("E51").Value = "  +2.81%  "
